# Insert a new weekly price record for "Jengibre" (ginger) at row 192,
# pushing the existing row 192 (and everything below it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(192).Insert()

$ws.Cells.Item(192, 1).Value = 10
$ws.Cells.Item(192, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(192, 3).Value = "La Araucanía"
$ws.Cells.Item(192, 4).Value = 45027
$ws.Cells.Item(192, 5).Value = 9
$ws.Cells.Item(192, 6).Value = 100114007
$ws.Cells.Item(192, 7).Value = "Jengibre"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Primera"
$ws.Cells.Item(192, 10).Value = 100
$ws.Cells.Item(192, 11).Value = 25000
$ws.Cells.Item(192, 12).Value = 25000
$ws.Cells.Item(192, 13).Value = 25000
$ws.Cells.Item(192, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(192, 15).Value = "Perú"
$ws.Cells.Item(192, 16).Value = 1923
$ws.Cells.Item(192, 17).Value = 13
$ws.Cells.Item(192, 18).Value = "Hortaliza"
